# Update "想去人数" (want-to-go count) figures in the F column across the
# four sheets of the workbook, reflecting refreshed data for the generated
# gh-pages output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 583
$ws.Range("F3").Value = 263
$ws.Range("F5").Value = 744
$ws.Range("F6").Value = 384
$ws.Range("F8").Value = 161
$ws.Range("F10").Value = 228
$ws.Range("F11").Value = 6074
$ws.Range("F14").Value = 501
$ws.Range("F24").Value = 319
$ws.Range("F25").Value = 1023
$ws.Range("F27").Value = 1836
$ws.Range("F28").Value = 497

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 270

# Sheet "本地生活" (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 259

# Sheet "全部类型" (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 259
$ws.Range("F3").Value = 583
$ws.Range("F4").Value = 263
$ws.Range("F6").Value = 744
$ws.Range("F8").Value = 384
$ws.Range("F10").Value = 161
$ws.Range("F12").Value = 228
$ws.Range("F13").Value = 6074
$ws.Range("F16").Value = 270
$ws.Range("F17").Value = 501
$ws.Range("F34").Value = 319
$ws.Range("F35").Value = 1023
$ws.Range("F37").Value = 1836
$ws.Range("F38").Value = 497

$wb.Save()
